# Applies the edit described by the diff:
#  - rename the worksheet from "o554F-HW40.xpc" to "o554F"
#  - append a new data row (row 16) to the sheet, copying the formatting
#    of row 15's label column (A) and reusing the existing shared string
#    "HexGrid-60degTilt5degRes" for column B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab.
$ws.Name = "o554F"

# Duplicate the (bold + bordered) formatting used by the row-index column
# on the preceding row, so the new cell reuses the existing style record
# instead of minting a new one.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

# Fill in the new row's values.
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 1.002839692827024
$ws.Range("D16").Value = 0.9662521697578933
$ws.Range("E16").Value = 1.006223748755438
$ws.Range("F16").Value = 1.002839692827024
$ws.Range("G16").Value = 0.9781768555407804
$ws.Range("H16").Value = 1.015563277946242
$ws.Range("I16").Value = 1.000406761700286
$ws.Range("J16").Value = 0.9662521697578933
$ws.Range("K16").Value = 0.9862379592566659
$ws.Range("L16").Value = 0.9945388260418451
$ws.Range("M16").Value = 0.9949104177546108
